$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-16 column A: relink to new variable id
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 1).Value = "pop_region__population_totale"
}

# Add new rows 17-20 for pop_region__type_region
$ws.Cells.Item(17, 1).Value = "pop_region__type_region"
$ws.Cells.Item(17, 2).Value = "urbaine"
$ws.Cells.Item(17, 3).Value = 9432

$ws.Cells.Item(18, 1).Value = "pop_region__type_region"
$ws.Cells.Item(18, 2).Value = "périurbaine"
$ws.Cells.Item(18, 3).Value = 43

$ws.Cells.Item(19, 1).Value = "pop_region__type_region"
$ws.Cells.Item(19, 2).Value = "rurale"
$ws.Cells.Item(19, 3).Value = 3434

$ws.Cells.Item(20, 1).Value = "pop_region__type_region"
$ws.Cells.Item(20, 2).Value = "montagne"
$ws.Cells.Item(20, 3).Value = 9481

# Resize table to include new rows
$tbl = $ws.ListObjects.Item("Tableau3")
$tbl.Resize($ws.Range("A1:C20"))

$ws.Columns.Item(1).ColumnWidth = 25.6640625
$ws.Columns.Item(2).ColumnWidth = 9.83203125

$ws.Range("A22").Select() | Out-Null
